$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (I1, J1) - match the same formatting as existing header cells (H1, etc.)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for columns I (I0) and J (IF), rows 2-17
$data = @(
    @(4, 4),
    @(5, 6),
    @(5, 6),
    @(4, 6),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(8, 9),
    @(4, 4),
    @(7, 7),
    @(4, 4),
    @(7, 8),
    @(8, 9),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
